$wb = $excel.ActiveWorkbook

# --- "Planning and Execution" sheet: shift the WWnn week labels down by one
#     (WW02->WW01, WW03->WW02, WW04->WW03, WW05->WW04) in the Start Date /
#     End Date columns (E, F). ---
$ws = $wb.Worksheets.Item("Planning and Execution")

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = "WW01"
    $ws.Cells.Item($r, 6).Value = "WW02"
}

for ($r = 11; $r -le 14; $r++) {
    $ws.Cells.Item($r, 5).Value = "WW03"
    $ws.Cells.Item($r, 6).Value = "WW03"
}

$ws.Cells.Item(15, 5).Value = "WW04"
$ws.Cells.Item(15, 6).Value = "WW04"

# --- Active sheet / selection: move from "QOR" to "Planning and Execution",
#     selecting cell F16 there. ---
$ws.Activate()
$ws.Range("F16").Select()
